$wb = $excel.ActiveWorkbook

# --- Sheet 1: "ランサーズ" -- insert a new newest row at the top of the data (row 2) ---
$ws1 = $wb.Worksheets.Item("ランサーズ")

# Drop existing hyperlink objects on the URL column before we shuffle rows around,
# otherwise their `ref` stays pinned to the old cell instead of following the data.
$ws1.Range("F2:F26").Hyperlinks.Delete()

# Shift all existing data rows (2-26) down by one row (3-27) by copying the values.
$srcRange = $ws1.Range("A2:H26")
$vals = $srcRange.Value2
$dstRange = $ws1.Range("A3:H27")
$dstRange.Value2 = $vals

# Clear row 2 completely so stray values (e.g. column H) from the old row don't linger.
$ws1.Range("A2:H2").Value = $null

# Write the newly scraped listing into row 2.
$ws1.Range("A2").Value = "2025-08-29 01:16:05"
$ws1.Range("B2").Value = "【急募】既存スプレッドシートを新アカウントに移動依頼"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5381748"
$ws1.Range("G2").Value = 13

# Rebuild the hyperlinks for every URL cell (F2:F27) from the (now correctly shifted) text.
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws1.Cells.Item($r, 6)
    $url = $cell.Value2
    $ws1.Hyperlinks.Add($cell, $url)
}

# --- Sheet 2: "統計" -- append the matching stats-log row ---
$ws2 = $wb.Worksheets.Item("統計")

$ws2.Range("A15").Value = "2025-08-29T01:16:05.746955"
$ws2.Range("B15").Value = 14
$ws2.Range("C15").Value = "全案件リスト"
$ws2.Range("D15").Value = 50
$ws2.Range("E15").Value = 7
$ws2.Range("F15").Value = 4
$ws2.Range("G15").Value = 14
